$wb = $excel.ActiveWorkbook

# Update "F" column (想去人数 / number of people interested) values
# across the four worksheets, per the source diff.

# --- Worksheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 3314
$ws.Range("F4").Value = 26
$ws.Range("F5").Value = 1365
$ws.Range("F9").Value = 195
$ws.Range("F10").Value = 54
$ws.Range("F11").Value = 8531
$ws.Range("F12").Value = 8531
$ws.Range("F13").Value = 470
$ws.Range("F15").Value = 132
$ws.Range("F16").Value = 98
$ws.Range("F17").Value = 322
$ws.Range("F18").Value = 130
$ws.Range("F19").Value = 78
$ws.Range("F20").Value = 7
$ws.Range("F21").Value = 351
$ws.Range("F22").Value = 10748
$ws.Range("F23").Value = 10748
$ws.Range("F24").Value = 285
$ws.Range("F35").Value = 2074
$ws.Range("F37").Value = 40
$ws.Range("F38").Value = 2114
$ws.Range("F39").Value = 888
$ws.Range("F40").Value = 4068
$ws.Range("F41").Value = 1626
$ws.Range("F43").Value = 2575
$ws.Range("F44").Value = 3020
$ws.Range("F45").Value = 1233
$ws.Range("F46").Value = 167
$ws.Range("F48").Value = 335
$ws.Range("F49").Value = 302
$ws.Range("F50").Value = 40
$ws.Range("F51").Value = 116

# --- Worksheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 5
$ws.Range("F8").Value = 7
$ws.Range("F9").Value = 2
$ws.Range("F11").Value = 31
$ws.Range("F14").Value = 2
$ws.Range("F16").Value = 9
$ws.Range("F21").Value = 28

# --- Worksheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 6
$ws.Range("F3").Value = 18

# --- Worksheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 3314
$ws.Range("F7").Value = 26
$ws.Range("F8").Value = 1365
$ws.Range("F13").Value = 7
$ws.Range("F15").Value = 195
$ws.Range("F16").Value = 8531
$ws.Range("F17").Value = 470
$ws.Range("F19").Value = 132
$ws.Range("F20").Value = 98
$ws.Range("F21").Value = 322
$ws.Range("F22").Value = 130
$ws.Range("F23").Value = 78
$ws.Range("F24").Value = 7
$ws.Range("F25").Value = 10748
$ws.Range("F26").Value = 285
$ws.Range("F28").Value = 18
$ws.Range("F32").Value = 2
$ws.Range("F38").Value = 2074
$ws.Range("F40").Value = 40
$ws.Range("F41").Value = 2114
$ws.Range("F42").Value = 888
$ws.Range("F44").Value = 1629
$ws.Range("F45").Value = 3020
$ws.Range("F47").Value = 1233
$ws.Range("F48").Value = 335
$ws.Range("F49").Value = 302
$ws.Range("F50").Value = 40
$ws.Range("F51").Value = 116
